# Restored from revision of admin on 05/29/2020 07:43:07 AM.TEST Author: admin. Type: SAVE.
# The only functional change in this revision is the value stored in cell C10
# of the "Rules" sheet (the "From" value for rule R30), which changes from 18 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
